# ============================================================================
# "test creation nouvelle vm de prepro"
# Adds a new pre-production (PPR) reverse-proxy VM entry to the "vms" sheet,
# removes the now-superfluous blank template row, and appends the matching
# dropdown-list entries (ESX subnets) to the "params" sheet.
# ============================================================================

$wb = $excel.ActiveWorkbook

$wsVms    = $wb.Worksheets.Item("vms")
$wsParams = $wb.Worksheets.Item("params")

# ----------------------------------------------------------------------------
# 1) "vms" sheet: remove the blank template row (old row 3) so data rows
#    shift up by one.
# ----------------------------------------------------------------------------
$wsVms.Rows.Item(3).Delete()

# ----------------------------------------------------------------------------
# 2) "vms" sheet: fill in row 2 with the new VM's data (reverse proxy
#    externe, pre-production / "nut-dmz-02").
# ----------------------------------------------------------------------------
$wsVms.Range("B2").Value = "VSL-PPR-RPE-001"
$wsVms.Range("C2").Value = "VSL-PPR-RPE-001"
$wsVms.Range("D2").Value = "Reverse Proxy externe"
$wsVms.Range("D2").NumberFormat = "[$-F400]h:mm:ss\ AM/PM"
$wsVms.Range("E2").Value = "vsl-ppr-rpe-001"
$wsVms.Range("F2").Value = "NUT-DMZ-DC02"
$wsVms.Range("G2").Value = "nut-dmz-02"
$wsVms.Range("H2").Value = "esx_lib2_item"
$wsVms.Range("J2").Value = "PPR_REVERSE_PROXY_EXT"
$wsVms.Range("L2").Value = "NUT_DMZ_EXT_DC2_01"
$wsVms.Range("N2").Value = 4096
$wsVms.Range("O2").Value = 2
$wsVms.Range("R2").Value = "172.24.74.1"
$wsVms.Range("T2").Value = "172.24.74.254"
$wsVms.Range("U2").Value = "RECETTE"

Write-Host "vms row updated"
